# Replace the placeholder stimulus set (3-digit PNGs, old word list) with
# the working set of cue sequences: 2-digit JPGs + new German verb list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, image, word, category
$rows = @(
    @(2, "face/face030.jpg", "gelten", "face"),
    @(3, "flower/flower016.jpg", "hauen", "flower"),
    @(4, "face/face003.jpg", "laufen", "face"),
    @(5, "flower/flower030.jpg", "sondern", "flower"),
    @(6, "face/face023.jpg", "schenken", "face"),
    @(7, "flower/flower004.jpg", "saufen", "flower"),
    @(8, "face/face029.jpg", "klappen", "face"),
    @(9, "flower/flower017.jpg", "fühlen", "flower"),
    @(10, "flower/flower012.jpg", "fesseln", "flower"),
    @(11, "face/face026.jpg", "bitten", "face"),
    @(12, "face/face006.jpg", "tagen", "face"),
    @(13, "face/face013.jpg", "loben", "face"),
    @(14, "flower/flower020.jpg", "füttern", "flower"),
    @(15, "face/face031.jpg", "starten", "face"),
    @(16, "flower/flower002.jpg", "währen", "flower"),
    @(17, "face/face010.jpg", "ehren", "face"),
    @(18, "flower/flower019.jpg", "drohen", "flower"),
    @(19, "face/face014.jpg", "schicken", "face"),
    @(20, "flower/flower011.jpg", "runden", "flower"),
    @(21, "flower/flower022.jpg", "sieben", "flower"),
    @(22, "flower/flower007.jpg", "liefern", "flower"),
    @(23, "flower/flower000.jpg", "langen", "flower"),
    @(24, "face/face025.jpg", "tauschen", "face"),
    @(25, "face/face011.jpg", "bleiben", "face"),
    @(26, "flower/flower025.jpg", "spielen", "flower"),
    @(27, "flower/flower015.jpg", "schmecken", "flower"),
    @(28, "face/face002.jpg", "biegen", "face"),
    @(29, "face/face007.jpg", "raten", "face"),
    @(30, "face/face027.jpg", "kehren", "face"),
    @(31, "face/face008.jpg", "hupen", "face"),
    @(32, "flower/flower026.jpg", "füllen", "flower"),
    @(33, "flower/flower013.jpg", "fliehen", "flower")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]   # B: image
    $ws.Cells.Item($r, 3).Value = $row[2]   # C: word
    $ws.Cells.Item($r, 4).Value = $row[3]   # D: category
}

